$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2636.1428
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 3045.111
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 3045.111
$ws.Range("M62").Value = -1276
$ws.Range("N62").Value = -4293.111
$ws.Range("H65").Value = 2636.1428
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 3045.111
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 15225.555
$ws.Range("M65").Value = -6380
$ws.Range("N65").Value = -21465.555
$ws.Range("H132").Value = 420571.22
$ws.Range("I132").Value = 438769.97
$ws.Range("K132").Value = 1316309.91
$ws.Range("M132").Value = -1313779.91

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2820.3572
$ws.Range("I88").Value = 3626.6667
$ws.Range("J88").Value = 2600.4546
$ws.Range("K88").Value = 3626.6667
$ws.Range("L88").Value = 2600.4546
$ws.Range("M88").Value = -3220.6667
$ws.Range("N88").Value = -3412.4546
$ws.Range("H91").Value = 2820.3572
$ws.Range("I91").Value = 3626.6667
$ws.Range("J91").Value = 2600.4546
$ws.Range("K91").Value = 3626.6667
$ws.Range("L91").Value = 2600.4546
$ws.Range("M91").Value = -2222.6667
$ws.Range("N91").Value = -5408.4546
$ws.Range("H122").Value = 2496
$ws.Range("I122").Value = 2496
$ws.Range("K122").Value = 7488
$ws.Range("M122").Value = -5038

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2717.65
$ws.Range("I86").Value = 2495.9092
$ws.Range("J86").Value = 2988.6667
$ws.Range("K86").Value = 2495.9092
$ws.Range("L86").Value = 2988.6667
$ws.Range("M86").Value = -1372.9092
$ws.Range("N86").Value = -5234.6667
$ws.Range("H89").Value = 2717.65
$ws.Range("I89").Value = 2495.9092
$ws.Range("J89").Value = 2988.6667
$ws.Range("K89").Value = 12479.546
$ws.Range("L89").Value = 14943.3335
$ws.Range("M89").Value = -6863.546
$ws.Range("N89").Value = -26175.3335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1923.7858
$ws.Range("J16").Value = 490
$ws.Range("L16").Value = 490
$ws.Range("N16").Value = -1064
$ws.Range("H31").Value = 1607
$ws.Range("I31").Value = 1040.85
$ws.Range("J31").Value = 2636.3635
$ws.Range("K31").Value = 1040.85
$ws.Range("L31").Value = 2636.3635
$ws.Range("M31").Value = -745.8499999999999
$ws.Range("N31").Value = -3226.3635
$ws.Range("H34").Value = 1607
$ws.Range("I34").Value = 1040.85
$ws.Range("J34").Value = 2636.3635
$ws.Range("K34").Value = 1040.85
$ws.Range("L34").Value = 2636.3635
$ws.Range("M34").Value = -838.8499999999999
$ws.Range("N34").Value = -3040.3635
$ws.Range("H62").Value = 3898.625
$ws.Range("I62").Value = 3838
$ws.Range("J62").Value = 3999.6667
$ws.Range("K62").Value = 3838
$ws.Range("L62").Value = 3999.6667
$ws.Range("M62").Value = -3214
$ws.Range("N62").Value = -5247.6667
$ws.Range("H65").Value = 3898.625
$ws.Range("I65").Value = 3838
$ws.Range("J65").Value = 3999.6667
$ws.Range("K65").Value = 19190
$ws.Range("L65").Value = 19998.3335
$ws.Range("M65").Value = -16070
$ws.Range("N65").Value = -26238.3335
$ws.Range("H99").Value = 1310.5
$ws.Range("I99").Value = 1292.4286
$ws.Range("J99").Value = 1328.5714
$ws.Range("K99").Value = 1292.4286
$ws.Range("L99").Value = 1328.5714
$ws.Range("M99").Value = 205.5714
$ws.Range("N99").Value = -4324.5714
$ws.Range("H107").Value = 376.13635
$ws.Range("I107").Value = 383.44446
$ws.Range("J107").Value = 343.25
$ws.Range("K107").Value = 383.44446
$ws.Range("L107").Value = 343.25
$ws.Range("M107").Value = 1536.55554
$ws.Range("N107").Value = -4183.25
$ws.Range("H113").Value = 1923.7858
$ws.Range("J113").Value = 490
$ws.Range("L113").Value = 490
$ws.Range("N113").Value = -4830
$ws.Range("H126").Value = 1310.5
$ws.Range("I126").Value = 1292.4286
$ws.Range("J126").Value = 1328.5714
$ws.Range("K126").Value = 3877.2858
$ws.Range("L126").Value = 3985.7142
$ws.Range("M126").Value = -1407.2858
$ws.Range("N126").Value = -8925.7142

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 413
$ws.Range("I36").Value = 101
$ws.Range("J36").Value = 725
$ws.Range("K36").Value = 303
$ws.Range("L36").Value = 2175
$ws.Range("M36").Value = -134
$ws.Range("N36").Value = -2513
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H49").Value = 2633.3333
$ws.Range("J49").Value = 2633.3333
$ws.Range("L49").Value = 7899.999899999999
$ws.Range("N49").Value = -8211.999899999999
$ws.Range("H110").Value = 3357.1428
$ws.Range("J110").Value = 3740
$ws.Range("L110").Value = 11220
$ws.Range("N110").Value = -19400
$ws.Range("H123").Value = 3000
$ws.Range("I123").Value = 3000
$ws.Range("K123").Value = 9000
$ws.Range("M123").Value = -6550
$ws.Range("H124").Value = 2180
$ws.Range("I124").Value = 950
$ws.Range("J124").Value = 3000
$ws.Range("K124").Value = 2850
$ws.Range("L124").Value = 9000
$ws.Range("M124").Value = 2060
$ws.Range("N124").Value = -18820
$ws.Range("H125").Value = 2200
$ws.Range("I125").Value = 500
$ws.Range("J125").Value = 3900
$ws.Range("K125").Value = 1500
$ws.Range("L125").Value = 11700
$ws.Range("M125").Value = 3420
$ws.Range("N125").Value = -21540
$ws.Range("H132").Value = 47619908
$ws.Range("I132").Value = 66667210
$ws.Range("K132").Value = 600004890
$ws.Range("M132").Value = -600002360
$ws.Range("H133").Value = 5537.294
$ws.Range("J133").Value = 7779.45
$ws.Range("L133").Value = 23338.35
$ws.Range("N133").Value = -33458.35

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2675.3845
$ws.Range("I132").Value = 2408.1428
$ws.Range("K132").Value = 7224.428400000001
$ws.Range("M132").Value = -4694.428400000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2363.4
$ws.Range("I7").Value = 2131.9092
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2131.9092
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2019.9092
$ws.Range("N7").Value = -3224
$ws.Range("H11").Value = 3000
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3280
$ws.Range("H126").Value = 2363.4
$ws.Range("I126").Value = 2131.9092
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6395.7276
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3925.7276
$ws.Range("N126").Value = -13940
$ws.Range("H136").Value = 2700.2727
$ws.Range("I136").Value = 1284
$ws.Range("J136").Value = 4399.8
$ws.Range("K136").Value = 3852
$ws.Range("L136").Value = 13199.4
$ws.Range("M136").Value = -1302
$ws.Range("N136").Value = -18299.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9312
$ws.Range("I136").Value = 17499.834
$ws.Range("J136").Value = 1124.1666
$ws.Range("K136").Value = 52499.50199999999
$ws.Range("L136").Value = 3372.4998
$ws.Range("M136").Value = -49949.50199999999
$ws.Range("N136").Value = -8472.4998

Write-Host "All changes applied."